# Daily attendance processing - 2026-01-25 17:34:14
# Swap the order of "dnasr281@gmail.com" and "System" in the
# "Recorded By" column (G) wherever both names are recorded together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
